# Update Name of Algo
# Apply updated KNN imputation results to specific cells on Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A4").Value = -21.913
$ws.Range("D4").Value = -8.132999999999999

$ws.Range("D5").Value = -8.616999999999999

$ws.Range("A6").Value = -21.14
$ws.Range("D6").Value = -8.309999999999999

$ws.Range("A7").Value = -21.018

$ws.Range("A8").Value = -20.727
$ws.Range("D8").Value = -8.334999999999999

$ws.Range("A16").Value = -20.727
$ws.Range("D16").Value = -8.405999999999999

$ws.Range("A20").Value = -21.86

$ws.Range("A21").Value = -21.14

$ws.Range("D22").Value = -8.16
